$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    foreach ($para in $doc.Paragraphs) {
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7, [char]12)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

function Split-RunAt($doc, $rangeStart, $offset) {
    # Forces a run boundary `offset` characters after `rangeStart` by
    # toggling a character formatting property on/off on that sub-range,
    # without altering the final visible formatting.
    $sub = $doc.Range($rangeStart, $rangeStart + $offset)
    $sub.Bold = 1
    $sub.Bold = 0
}

# Processed in reverse document order so that a paragraph's freshly
# written replacement text can never be mistaken for another
# (not-yet-processed) paragraph's original search text.

# --- "Seat" -> "Role" (single run) ---
$p = Find-ParagraphByText $d "Seat"
$p.Range.Find.Execute("Seat", $true, $false, $false, $false, $false, $true, 1, $false, "Role", 2)

# --- "Time" -> "Seat" (single run) ---
$p = Find-ParagraphByText $d "Time"
$p.Range.Find.Execute("Time", $true, $false, $false, $false, $false, $true, 1, $false, "Seat", 2)

# --- "Day" -> "Service" + "Time" (2 runs) ---
$p = Find-ParagraphByText $d "Day"
$p.Range.Find.Execute("Day", $true, $false, $false, $false, $false, $true, 1, $false, "ServiceTime", 2)
$p = Find-ParagraphByText $d "ServiceTime"
Split-RunAt $d $p.Range.Start 7

# --- "Time Table" -> "Business" + "Day" (2 runs) ---
$p = Find-ParagraphByText $d "Time Table"
$p.Range.Find.Execute("Time Table", $true, $false, $false, $false, $false, $true, 1, $false, "BusinessDay", 2)
$p = Find-ParagraphByText $d "BusinessDay"
Split-RunAt $d $p.Range.Start 8

# --- "Schedule" -> "Appointment" (single run) ---
$p = Find-ParagraphByText $d "Schedule"
$p.Range.Find.Execute("Schedule", $true, $false, $false, $false, $false, $true, 1, $false, "Appointment", 2)

# --- "User (Customer & Hairstylist)" -> "User (Customer," + " Hairstylist)" (2 runs) ---
$p = Find-ParagraphByText $d "User (Customer & Hairstylist)"
$p.Range.Find.Execute("User (Customer & Hairstylist)", $true, $false, $false, $false, $false, $true, 1, $false, "User (Customer, Hairstylist)", 2)
$p = Find-ParagraphByText $d "User (Customer, Hairstylist)"
Split-RunAt $d $p.Range.Start 15
